$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Password" header in column J. (Adding this right after freeing the
# old "sur123" shared-string slot below lets it land back in that slot.)
$ws.Range("J1").Value = "Password"

# Quote-wrapped literal text values replacing the old numeric entries.
$ws.Range("F2").Value = '"580024"'
$ws.Range("G2").Value = '"1234567891"'
$ws.Range("H2").Value = '"123"'
$ws.Range("I2").Value = '"sur123"'

# Password value (leading space, as typed by the author) + new Confirm column.
$ws.Range("J2").Value = " Suraj@123"
$ws.Range("K1").Value = "Confirm"
$ws.Range("K2").Value = " Suraj@123"

# Turn the password cell into a hyperlink (adds the Hyperlink style/font
# and the worksheet <hyperlinks> entry), then restore the leading-space
# cell text (Hyperlinks.Add overwrites the cell with the display text).
[void]$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:Suraj@123", "", "", "Suraj@123")
$ws.Range("J2").Value = " Suraj@123"

# Move the active selection to K2, matching the saved view state.
[void]$ws.Range("K2").Select()
